$wb = $excel.ActiveWorkbook

# Sheet 1: "VENTAS POR GRUPO"
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L27").Value = 346.56
$ws1.Range("L28").Value = "2 de 26"

# Sheet 2: "VENTA MENSUAL"
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F27").Value = 346.56
$ws2.Range("F28").Value = 8670.459999999999

# Sheet 3: "CUMPLIMIENTO MENSUAL"
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D2").Value = 4429.98
$ws3.Range("E2").Value = -4429.98
$ws3.Range("D4").Value = 14124.06
$ws3.Range("E4").Value = 5875.940000000001
$ws3.Range("F4").Value = 0.706203
